{"js": "// Each 'old' text below is unique in the document body, so a\n// matchCase search is guaranteed to resolve to exactly one run\n// and insertText(..., Word.InsertLocation.replace) swaps its\n// text while keeping the run's existing formatting.\nconst replacements = [\n  [\"2025-09-16 Tuesday\", \"2025-09-17 Wednesday\"],\n  [\"605\u00d74=2420\", \"854\u00d77=5978\"],\n  [\"573\u00d79=5157\", \"372\u00d76=2232\"],\n  [\"241\u00d79=2169\", \"342\u00d74=1368\"],\n  [\"729\u00d75=3645\", \"685\u00d75=3425\"],\n  [\"953\u00d77=6671\", \"595\u00d74=2380\"],\n  [\"364\u00d77=2548\", \"638\u00d74=2552\"],\n  [\"217\u00d78=1736\", \"101\u00d74=404\"],\n  [\"578\u00d72=1156\", \"238\u00d72=476\"],\n  [\"190\u00d76=1140\", \"915\u00d77=6405\"],\n  [\"952\u00d79=8568\", \"143\u00d78=1144\"],\n  [\"460\u00d73=1380\", \"626\u00d79=5634\"],\n  [\"786\u00d77=5502\", \"542\u00d74=2168\"],\n  [\"163\u00d78=1304\", \"433\u00d78=3464\"],\n  [\"556\u00d74=2224\", \"373\u00d77=2611\"],\n  [\"214\u00d74=856\", \"740\u00d78=5920\"],\n  [\"353\u00d76=2118\", \"157\u00d79=1413\"],\n  [\"769\u00d76=4614\", \"411\u00d73=1233\"],\n  [\"365\u00d78=2920\", \"208\u00d72=416\"],\n  [\"309\u00d76=1854\", \"708\u00d76=4248\"],\n  [\"660\u00d75=3300\", \"355\u00d76=2130\"],\n  [\"890\u00d76=5340\", \"849\u00d77=5943\"],\n  [\"554\u00d74=2216\", \"465\u00d73=1395\"],\n  [\"379\u00d77=2653\", \"457\u00d76=2742\"],\n  [\"498\u00d72=996\", \"684\u00d79=6156\"],\n  [\"112\u00d73=336\", \"691\u00d73=2073\"]\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each date/equation text run to the new value.\n# Every \"old\" string in the document is unique, so exact\n# Find/Replace (MatchCase, no wildcards) cannot cross-match.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2025-09-16 Tuesday'; New = '2025-09-17 Wednesday' },\n    @{ Old = '605\u00d74=2420'; New = '854\u00d77=5978' },\n    @{ Old = '573\u00d79=5157'; New = '372\u00d76=2232' },\n    @{ Old = '241\u00d79=2169'; New = '342\u00d74=1368' },\n    @{ Old = '729\u00d75=3645'; New = '685\u00d75=3425' },\n    @{ Old = '953\u00d77=6671'; New = '595\u00d74=2380' },\n    @{ Old = '364\u00d77=2548'; New = '638\u00d74=2552' },\n    @{ Old = '217\u00d78=1736'; New = '101\u00d74=404' },\n    @{ Old = '578\u00d72=1156'; New = '238\u00d72=476' },\n    @{ Old = '190\u00d76=1140'; New = '915\u00d77=6405' },\n    @{ Old = '952\u00d79=8568'; New = '143\u00d78=1144' },\n    @{ Old = '460\u00d73=1380'; New = '626\u00d79=5634' },\n    @{ Old = '786\u00d77=5502'; New = '542\u00d74=2168' },\n    @{ Old = '163\u00d78=1304'; New = '433\u00d78=3464' },\n    @{ Old = '556\u00d74=2224'; New = '373\u00d77=2611' },\n    @{ Old = '214\u00d74=856'; New = '740\u00d78=5920' },\n    @{ Old = '353\u00d76=2118'; New = '157\u00d79=1413' },\n    @{ Old = '769\u00d76=4614'; New = '411\u00d73=1233' },\n    @{ Old = '365\u00d78=2920'; New = '208\u00d72=416' },\n    @{ Old = '309\u00d76=1854'; New = '708\u00d76=4248' },\n    @{ Old = '660\u00d75=3300'; New = '355\u00d76=2130' },\n    @{ Old = '890\u00d76=5340'; New = '849\u00d77=5943' },\n    @{ Old = '554\u00d74=2216'; New = '465\u00d73=1395' },\n    @{ Old = '379\u00d77=2653'; New = '457\u00d76=2742' },\n    @{ Old = '498\u00d72=996'; New = '684\u00d79=6156' },\n    @{ Old = '112\u00d73=336'; New = '691\u00d73=2073' }\n)\n\nforeach ($pair in $replacements) {\n    $findRange = $d.Content\n    $findRange.Find.ClearFormatting()\n    $findRange.Find.Replacement.ClearFormatting()\n    # wdFindContinue=1, wdReplaceAll=2 -- replace every occurrence\n    # (there is exactly one per string) in a single Execute call.\n    $findRange.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
